# Populate 'sm_usoyforma' worksheet (sheet8.xml) with 28 new rows
# as part of creating the 'Uso y Forma' DB insert script.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sm_usoyforma")

$b1 = @'
Cántaros
'@
$c1 = @'
Las formas básicas proceden del Preclásico Temprano y Medio, donde el inventario cerámico incluye grandes tinajas sin cuello.
'@

$b2 = @'
Cuencos
'@
$c2 = @'
Según la concepción más aceptada así se le denomina a las vasijas cuyo diámetro es igual o mayor a su altura, con variantes en la profundidad y diámetro sin llegar al extremo como en el caso de los platos.
'@

$b3 = @'
Esculpidos / Tallados
'@
$c3 = @'
La destreza artística de los talladores y escultores mayas quedó manifiesta en diversidad de superficies. Fue amplia la cantidad de monumentos como estelas, dinteles y altares, pero también hubo preferencia por objetos de menor tamaño incluyendo artefactos portátiles. La colección aquí representada que incluye diversidad de usos y formas como fragmentos de altar, mascarones de estuco, placas de piedra verde, piedras de moler, yugos, así como talla en hueso y concha.
'@

$b4 = @'
Figurilla
'@
$c4 = @'
Las figurillas son objetos mayoritariamente elaborados en cerámica, aunque los hay en piedra verde e incluso en hueso como los aquí mostrados. Generalmente son representaciones antropomorfas que se han encontrado desde los estratos más antiguos junto a las primeras cerámicas, entonces fueron modeladas a mano. Las más antiguas figurillas presentan características individuales y no hay una homogeneidad de estilos, si bien fue una tradición ampliamente difundida en Mesoamérica. El inventario de este tipo de objetos se incrementó durante el Preclásico Medio y Tardío aunque comenzaron a ser más uniformes en su estilo. Entrado el Clásico Temprano hubo un cese en su fabricación, lo cual sigue siendo un asunto intrigante, pero finalmente la antigua tradición fue retomada pero con ya notables modificaciones durante el Clásico Tardío.
'@

$b5 = @'
Instrumentos Sonoros
'@
$c5 = @'
Una característica casi universal de todas las culturas es el apego hacia la creación musical, práctica tan antigua como el mismo ser humano y de la cual los mayas prehispánicos también supieron cultivar. No obstante las variaciones de instrumentos solamente contemplaron dos tipos en general: los del tipo membranófonos, como tambores y percutores; y los aerófonos, es decir los instrumentos de viento como trompetas, silbatos y flautas.\nLa colección que aquí se presenta contiene solamente del tipo aerófono, presentando variación en cuanto a material y forma, acentuando la preferencia por las formas humanas elaboradas en barro cocido, aunque las elaboradas en concha son destacadas y atractivas excepciones.
'@

$b6 = @'
Incensarios
'@
$c6 = @'
Es una vasija de diversas formas, que puede diferir del contexto o procedencia y sobre todo es variante en el tiempo. Su función básica es la de quemar incienso, lo cual le da una connotación ritual. Este tipo de vasijas está presente durante toda la ocupación prehispánica e incluso sobrevivió al proceso de conquista.
'@

$b7 = @'
Morteros
'@
$c7 = @'
Son utensilios, preferentemente de piedra, que sirven para machacar distintas sustancias como especias, semillas o extractas alucinógenas. Generalmente utilizados en la preparación de alimentos.
'@

$b8 = @'
Ollas
'@
$c8 = @'
Las ollas son cuencos profundos y con boca ancha, usualmente es un artefacto de cocina que servía para hervir algún alimento.
'@

$b9 = @'
Piedras de moler
'@
$c9 = @'
Como su nombre lo indica es un instrumento que sirve para la molienda de semillas y usualmente es cóncavo. Este tipo de artefacto fue ampliamente difundido en toda Mesoámerica, pues conlleva una actividad cotidiana y de sustento. Tales instrumentos suelen estar acompañados de las manos de moler y a este conjunto se le conoce como metate. Una característica en el desarrollo histórico de las piedras de moler es que las fabricadas durante el Preclásico carecen de soportes, mientras que las elaboradas en el periodo Clásico pueden ser trípodes o tetrápodos.
'@

$b10 = @'
Sellos
'@
$c10 = @'
Los sellos representan un elemento común en diversos contextos arqueológicos a lo largo del área mesoamericana, presentes prácticamente en toda la época prehispánica. Estos fueron principalmente elaborados en barro, si bien existen algunos en piedra. Estos fueron planos o cilíndricos, donde los primeros se estampaban mediante la presión y los segundos mediante la rotación del sello sobre una superficie. Muchas representaciones iconográficas muestran grabados en la piel de diversos personajes que se corresponden con formas encontradas en diversos sellos.\nEl ejemplar que aquí se presenta contiene motivos zoomorfos, posiblemente sean monos con integración de elementos fitomorfos, es decir elementos vegetales.
'@

$b11 = @'
Urnas
'@
$c11 = @'
Por lo general cuencos grandes y profundos que tienen varios usos pero primordialmente para guardar los restos o cenizas de difuntos. Muchas de las urnas portan motivos antropomorfos y zoomorfos o una combinación de ambas incluyendo soportes.
'@

$b12 = @'
Vasijas
'@
$c12 = @'
Las vasijas de tipo zapato son características del Preclásico Medio, también las vasijas con vertedera.
'@

$b13 = @'
Vasos
'@
$c13 = @'
Los vasos en un nivel muy general pueden denominarse como vasijas cuya altura siempre es mayor que su diámetro. Por lo general las paredes de los vasos son verticales, aunque los hay con leves curvas divergentes o convergentes, las bases pueden ser planas o con soportes de tipo tetrápode o trípode. En cuanto a las decoraciones estas son muy variantes destacando la policromía durante el Clásico o los denominados vasos tipo códice con representación de escenas míticas y escritura jeroglífica.
'@

$b14 = @'
Glíficos
'@
$c14 = @'
Se dice que uno de los logros más notables en la historia del ser humano fue la invención de la escritura. Avance científico e intelectual que no muchas culturas de la antigüedad lograron desarrollar, pero que la cultura maya puede presumir estar dentro de ese selecto número de sociedades capaces de crear un sistema que expresara el lenguaje hablado en caracteres escritos. Este sistema de escritura comúnmente llamado jeroglífico, en referencia al sistema creado en el antiguo Egipto, fue un complejo conjunto de símbolos desarrollado posiblemente desde el Preclásico Tardío y perfeccionado durante el Clásico. Conformado por aproximadamente 800 símbolos diferentes que integran un sistema mixto de tipo logo-silábico, es decir, que emplea logogramas que significan palabras completas, sílabas y vocales que complementan el sistema.\nFue tan amplia la difusión de este sistema entre los grupos de poder de las ciudades mayas que pueden encontrarse referencias escritas en casi todos los materiales conocidos de aquella época. Generalmente la escritura fue desarrollada como medio propagandístico de las élites gobernantes, principalmente en los grandes monumentos que relatan la historia política y dinástica de las ciudades clásicas. Por otro lado los objetos portátiles con escritura hacen referencia al uso de tales objetos, nombres de artistas, escribanos o incluso algunas narraciones míticas.\nEn la muestra de esta aplicación se presentan algunos ejemplares de cerámica con escritura maya, además de un fragmento de altar y el recién descubierto Panel V de la Corona, el cual hace referencia a la muy promovida fecha 4 ajaw 3 kank’in, correspondiente al 21 de diciembre de 2012.
'@

$b15 = @'
Jars
'@
$c15 = @'
The most basic forms date from the Early and Main Pre-Classical Period, in which the inventory of ceramics consists primarily of big clay jugs without bottleneck.
'@

$b16 = @'
Bowls
'@
$c16 = @'
According to the most common definition, a bowl is a type of container, that has a diameter equal or mayor compared to its height. There are variations concerning profoundness and diameter, without ever reaching an extreme relation like plates.
'@

$b17 = @'
Gouging/Engraving
'@
$c17 = @'
The artistic craftsmanship of Mayan engravers and artisans skilled in gouging techniques has been conserved in the diversity of surfaces. There was a large quantity of monuments like steles, lintels and altars, even if there was a significant preference for smaller artefacts, like portable objects as well. The present collection includes distinct uses and forms like altar fragments, figureheads made of stucco, boards made of green stone, millstones, yokes, and engraved shells and bones.
'@

$b18 = @'
Figurines
'@
$c18 = @'
Figurines are objects mainly fabricated of ceramic material, although there are samples made from green stone and even bone, like the examples shown here. They are in general anthropomorphic representations that have been found since the eldest strata at the same time as the first ceramic works were found. The Figurines were moulded in handicraft. The most antique Figurines represent individual characteristics with no stylistic homogeneity. This tradition was nevertheless widely spread in Mesoamerica. The inventory of this type of objects was incremented during the Main and Late Pre-Classical Period, though their style became more uniform. By the start of the Classical Period their fabrication ceased, which is an interesting issue. The tradition was however reanimated during the Late Classical Period, but with notable modifications.
'@

$b19 = @'
Audio-Instruments
'@
$c19 = @'
An almost universal characteristic of every culture is the adherence to the creation of music, a practise as old as humanity itself and which the Maya also cultivated. Yet the variation of musical instruments includes only two different types: the membrane-based type, like drums and percussion instruments, aero-phone instruments, wind-instruments so to speak, such as trumpets, whistles and flutes.\nThe here presented collection contains exclusively instruments of the aero-phone type, presenting variations concerning material and form, in which a clear preference for humanly constructed forms (terracotta) is expressed, even though the ones constructed from shells are an impressive and outstanding exception.
'@

$b20 = @'
Censers (Thuribles)
'@
$c20 = @'
These are diversely shaped containers, differing after their context and origin (but in the first place over time), that are basically used for the religious cremation of incense, which gives the object a ritualistic connotation. This type of containers was present during the whole Pre-Hispanic occupation and even survived the process of the Conquista.
'@

$b21 = @'
Mortars
'@
$c21 = @'
These are objects, of stone by preference, that are used to mash different materials like spices, seeds or hallucinogenic extracts. They are principally used for the preparation of meals.
'@

$b22 = @'
Pots
'@
$c22 = @'
Pots are deep bowls with wide mouths, they are usually kitchen artefacts that served in general to boil water and cook food.
'@

$b23 = @'
Millstones
'@
$c23 = @'
According to their name, they are instruments used for the milling of seeds and usually they are concave. This type of artefact was widely spread in all Mesoamerica, since it entails an everyday activity that contributes to the means of subsistence. These instruments are usually accompanied by a certain type of hand-axe. This set of items is known by the word metate. A characteristic aspect in the development of millstones is the fact that the examples dating from the Pre-Classical Period dispense with a frame, while the more elaborate ones from the Classical Period tend to be tripods or even four-footed.
'@

$b24 = @'
Seals
'@
$c24 = @'
The seals represent a common element in divers archaeological contexts throughout the Mesoamerican area, since they are present in practically all the Pre-Hispanic era. They were primordially worked in clay, although there are some examples made of stone. They were flat or cylindrical, of which the first kind was stamped applying pressure on top, and the second kind by rolling the seal over a surface. Many iconographic representations show engravings in the skin of different persons that correspond to some forms found upon distinct seals.\nThe here featured example contains animal related themes, possibly apes combined with motifs of vegetation.
'@

$b25 = @'
Urns
'@
$c25 = @'
Urns are in general big and profound bowls that have distinct functions, but which are primarily used to preserve the rests or ashes of the deceased. Many of the urns feature anthropomorphic motifs and iconographic representations of animals or a combination of both. Even racks for urns have been found.
'@

$b26 = @'
Containers
'@
$c26 = @'
The shoe shaped containers are a characteristic aspect of the Main Pre-Classical Period, as well as the ones with mouldboards.
'@

$b27 = @'
Cups
'@
$c27 = @'
Cups can be defined – in a very general way – as containers whose height is always mayor in comparison to their diameter. The cups""' walls are generally vertical, although some examples show slight divergent or convergent curves. Their bases can be both flat or with tripod support frames. The decorations of these cups are very variable, and it""'s their polychromatic design that is outstanding during the Classical Period. Some cups belong to rituals and carry iconographic representations of mythological scenes and hieroglyphic scripture.
'@

$b28 = @'
Glyphs
'@
$c28 = @'
It is said that one of the most notable merits in the history of mankind is the invention of writing. It is considered an intellectual and scientific progress that not many cultures in the antiquity achieved to develop. The Mayan culture however can pride itself of appertaining to the select number of societies that has been able to create a system of expressing spoken language in written characters. This writing system commonly known as Hieroglyphs, referring to the system created in ancient Egypt, was a complex conjunct of symbols, probably developed since the Late Pre-Classical Period and perfected during the Classical Period. Composed by approximately 800 different symbols the Mayan writing is a mixed, logo-syllabic writing system, this is to say that it uses logograms (symbols that designate complete words), syllables, and vowels that complete the system.\nThe diffusion of this system among the different groups of power in the Mayan cities was so vast that written references can be found upon nearly every material known at the time. The scripture was generally developed to the means of propaganda of the ruling elites, principally upon the great monuments that tell the political and dynastic history of the Classical cities. On the other hand there are portable objects, that have been written on, which refer to the use of the very object, names of artists, writers or even some mythic narrations.\nAmong the here presented samples there are some ceramic samples carrying Mayan scripture upon. There is also a fragment of an altar and the recently discovered Panel V of La Corona, which alludes to the very promoted date 4 Ajaw 3 Kank""'in, that corresponds to December 21. in 2012.
'@

$bArr = @($b1,$b2,$b3,$b4,$b5,$b6,$b7,$b8,$b9,$b10,$b11,$b12,$b13,$b14,$b15,$b16,$b17,$b18,$b19,$b20,$b21,$b22,$b23,$b24,$b25,$b26,$b27,$b28)
$cArr = @($c1,$c2,$c3,$c4,$c5,$c6,$c7,$c8,$c9,$c10,$c11,$c12,$c13,$c14,$c15,$c16,$c17,$c18,$c19,$c20,$c21,$c22,$c23,$c24,$c25,$c26,$c27,$c28)
$dArr = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,2,2,2,2,2,2,2,2,2,2,2,2,2,2)

for ($i = 0; $i -lt 28; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = ($i + 1)
    $ws.Cells.Item($row, 2).Value = $bArr[$i]
    $ws.Cells.Item($row, 3).Value = $cArr[$i]
    $ws.Cells.Item($row, 4).Value = $dArr[$i]

    $f = 'CONCATENATE("INSERT INTO sm_usoyforma VALUES (",A' + $row + ',",","""",B' + $row + ',"""",",","""",C' + $row + ',"""",",",D' + $row + ',");")'
    $ws.Cells.Item($row, 5).Formula = "=" + $f
}

$ws.Activate()
$ws.Range("A2:E29").Select()

Write-Host "Done populating sm_usoyforma with 28 rows"